$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New translation rows for the income chart legend keys (highest, lowest,
# daily, hourly, weekly) x 10 locales each. These 50 rows are inserted right
# before the old "language name" legend rows, pushing those down by 50.
$data = @(
  @(1, "highest", "Highest"),
  @(2, "highest", "மிக உயர்ந்தது"),
  @(3, "highest", "le plus élevé"),
  @(4, "highest", "最高"),
  @(5, "highest", "最高"),
  @(6, "highest", "Altíssima"),
  @(7, "highest", "उच्चतम"),
  @(8, "highest", "наибольший"),
  @(9, "highest", "Mas alto"),
  @(10, "highest", "الأعلى"),
  @(1, "lowest", "Lowest"),
  @(2, "lowest", "குறைந்த"),
  @(3, "lowest", "le plus bas"),
  @(4, "lowest", "最低"),
  @(5, "lowest", "最低"),
  @(6, "lowest", "mais baixo"),
  @(7, "lowest", "सबसे कम"),
  @(8, "lowest", "самый низкий"),
  @(9, "lowest", "El mas bajo"),
  @(10, "lowest", "أدنى"),
  @(1, "daily", "Daily"),
  @(2, "daily", "தினசரி"),
  @(3, "daily", "du quotidien"),
  @(4, "daily", "日常的"),
  @(5, "daily", "毎日"),
  @(6, "daily", "diário"),
  @(7, "daily", "रोज"),
  @(8, "daily", "повседневная"),
  @(9, "daily", "Diario"),
  @(10, "daily", "اليومي"),
  @(1, "hourly", "Hourly"),
  @(2, "hourly", "மணிக்கு "),
  @(3, "hourly", "toutes les heures"),
  @(4, "hourly", "每小时"),
  @(5, "hourly", "毎時"),
  @(6, "hourly", "de hora em hora"),
  @(7, "hourly", "प्रति घंटा"),
  @(8, "hourly", "ежечасно"),
  @(9, "hourly", "cada hora"),
  @(10, "hourly", "ساعيا"),
  @(1, "weekly", "Weekly"),
  @(2, "weekly", "வாரந்தோறும்"),
  @(3, "weekly", "hebdomadaire"),
  @(4, "weekly", "每周"),
  @(5, "weekly", "毎週"),
  @(6, "weekly", "semanalmente"),
  @(7, "weekly", "साप्ताहिक"),
  @(8, "weekly", "еженедельно"),
  @(9, "weekly", "semanalmente"),
  @(10, "weekly", "أسبوعي"),
)

$firstRow = 3792
$insertCount = $data.Count
$lastNewRow = $firstRow + $insertCount - 1

# Insert blank rows before the footer rows; this shifts rows 3792:3801 down
# to 3842:3851 and (since Excel inherits formatting from the row above on
# insert) gives the new B/D cells the same styles as row 3791 automatically.
$ws.Range("A" + $firstRow + ":A" + $lastNewRow).EntireRow.Insert() | Out-Null

for ($i = 0; $i -lt $insertCount; $i++) {
    $row = $firstRow + $i
    $prevRow = $row - 1
    $locId = $data[$i][0]
    $key = $data[$i][1]
    $text = $data[$i][2]

    $ws.Range("A" + $row).Formula = "=A" + $prevRow + "+1"
    $ws.Range("B" + $row).Value = $locId
    $ws.Range("C" + $row).Value = $key
    $ws.Range("D" + $row).Value = $text
    $ws.Range("E" + $row).Formula = "=CONCATENATE(" + [char]34 + "(" + [char]34 + ",CHAR(34),A" + $row + ",CHAR(34)," + [char]34 + "," + [char]34 + ",CHAR(34),B" + $row + ",CHAR(34)," + [char]34 + "," + [char]34 + ",CHAR(34),C" + $row + ",CHAR(34)," + [char]34 + "," + [char]34 + ",CHAR(34),D" + $row + ",CHAR(34)," + [char]34 + ")," + [char]34 + ")"
}

# Row 3821 ("daily" / locale 10 / Arabic) keeps the default (unstyled) D
# cell in the source workbook instead of inheriting the wrap-text style -
# match that quirk exactly.
$ws.Range("D3821").Style = "Normal"

# Restore the selection shown in the edited workbook (now pointing at the
# shifted legend row).
$ws.Range("E3846").Select() | Out-Null

"ok"
